$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.206.23'
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").Value = '2.174.88'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.36%  '
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '70.26'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.25%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.580'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.13'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0928'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.80'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.20%  '
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").Value = '2.500.59'
$ws.Range("E14").Value = '  -2.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.809'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.52%  '
$ws.Range("D17").Value = '2.160.18'
$ws.Range("E17").Value = '  -2.46%  '
$ws.Range("D18").Value = '41.029.27'
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("E19").Value = '  -7.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.87%  '
$ws.Range("E23").Value = '  -2.00%  '
$ws.Range("E24").Value = '  -7.10%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.22%  '
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.19%  '
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("E31").Value = '  -3.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0770'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  -9.39%  '
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("E36").Value = '  -8.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.13'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0285'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.19'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.55%  '
$ws.Range("E40").Value = '  -2.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.44'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '60.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.91%  '
$ws.Range("E43").Value = '  -4.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.97%  '
$ws.Range("E45").Value = '  -3.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '98.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.82%  '
$ws.Range("E47").Value = '  -2.91%  '
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.23'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.78%  '
$ws.Range("E50").Value = '  -3.00%  '
$ws.Range("D51").Value = '2.379.32'
$ws.Range("E51").Value = '  -2.00%  '
